$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert 4 new rows before the final "note" row (old row 33) to make room for
# a new "begin screen" / questions / "end screen" block.
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# Row 33: begin screen
$ws.Cells.Item(33, 1).Value = "begin screen"

# Row 34: Total Hours of Operation Unit 1 (integer question)
$ws.Cells.Item(34, 3).Value = "integer"
$ws.Cells.Item(34, 5).Value = "total_hrs_of_operation_unit1"

# Row 35: Total Hours of Operation Unit 2 (integer question)
$ws.Cells.Item(35, 3).Value = "integer"
$ws.Cells.Item(35, 5).Value = "total_hrs_of_operation_unit2"

# Titles (English) for both rows
$ws.Cells.Item(34, 6).Value = "Total Hours of Operation Unit 1"
$ws.Cells.Item(35, 6).Value = "Total Hours of Operation Unit 2"

# Titles (Spanish) for both rows
$ws.Cells.Item(34, 7).Value = "Total de horas de operación unidad 1"
$ws.Cells.Item(35, 7).Value = "Total de horas de operación unidad 2"

# Prompts (English) - unit 2 entered before unit 1
$ws.Cells.Item(35, 9).Value = "Enter total hours of operation for unit 2"
$ws.Cells.Item(34, 9).Value = "Enter total hours of operation for unit 1"

# Prompts (Spanish)
$ws.Cells.Item(34, 10).Value = "Ingrese el total de horas de operación para la unidad 1"
$ws.Cells.Item(35, 10).Value = "Ingrese el total de horas de operación para la unidad 2"

# Highlight the two new question rows (yellow fill, same as other newly-added rows)
$cols = @(3, 5, 6, 7, 9, 10)
foreach ($col in $cols) {
    $ws.Cells.Item(34, $col).Interior.Color = 65535
    $ws.Cells.Item(35, $col).Interior.Color = 65535
}

# Row 36: end screen
$ws.Cells.Item(36, 1).Value = "end screen"

# Update the sheet view to match the saved selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("N34").Select()
